$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D, shifting D:K -> E:L
$ws.Range("D1").EntireColumn.Insert()

# Copy number formats from the (now-shifted) neighbor column E so the new
# column D cells carry the same style as the data they sit beside (Excel
# does not auto-propagate formatting across an inserted column here).
$ws.Range("E7:E35").Copy()
$ws.Range("D7:D35").PasteSpecial(-4122)
$ws.Range("E38:E77").Copy()
$ws.Range("D38:D77").PasteSpecial(-4122)
$ws.Range("E80:E102").Copy()
$ws.Range("D80:D102").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Match the new columns width to its neighbors
$ws.Columns("D").ColumnWidth = $ws.Columns("E").ColumnWidth

# Populate the new column D with the latest reporting periods figures
$ws.Range("D7").Value2 = 43465
$ws.Range("D8").Value2 = 3962000
$ws.Range("D9").Value2 = 2661000
$ws.Range("D10").Value2 = 1301000
$ws.Range("D11").ClearContents()
$ws.Range("D12").Value2 = "NA"
$ws.Range("D13").Value2 = 0
$ws.Range("D14").Value2 = 0
$ws.Range("D15").Value2 = 556000
$ws.Range("D16").ClearContents()
$ws.Range("D17").Value2 = 3371000
$ws.Range("D18").Value2 = 591000
$ws.Range("D19").ClearContents()
$ws.Range("D20").Value2 = 17000
$ws.Range("D21").Value2 = 1164000
$ws.Range("D22").Value2 = 115000
$ws.Range("D23").Value2 = 493000
$ws.Range("D24").Value2 = 92000
$ws.Range("D25").Value2 = 0
$ws.Range("D26").Value2 = 401000
$ws.Range("D27").Value2 = 400000
$ws.Range("D28").Value2 = 0
$ws.Range("D29").Value2 = "NA"
$ws.Range("D30").Value2 = 0
$ws.Range("D31").Value2 = 0
$ws.Range("D32").Value2 = -17000
$ws.Range("D33").Value2 = 400000
$ws.Range("D34").Value2 = 0
$ws.Range("D35").Value2 = 400000
$ws.Range("D38").Value2 = 43465
$ws.Range("D39").ClearContents()
$ws.Range("D40").ClearContents()
$ws.Range("D41").Value2 = 18000
$ws.Range("D42").Value2 = 0
$ws.Range("D43").Value2 = 738000
$ws.Range("D44").Value2 = 134000
$ws.Range("D45").Value2 = 80000
$ws.Range("D46").Value2 = 970000
$ws.Range("D47").Value2 = 0
$ws.Range("D48").Value2 = 12439000
$ws.Range("D49").Value2 = 0
$ws.Range("D50").Value2 = 0
$ws.Range("D51").Value2 = 0
$ws.Range("D52").Value2 = 1980000
$ws.Range("D53").Value2 = 0
$ws.Range("D54").Value2 = 15389000
$ws.Range("D55").ClearContents()
$ws.Range("D56").ClearContents()
$ws.Range("D57").Value2 = 556000
$ws.Range("D58").Value2 = 259000
$ws.Range("D59").Value2 = 943000
$ws.Range("D60").Value2 = 1758000
$ws.Range("D61").Value2 = 3427000
$ws.Range("D62").Value2 = 5946000
$ws.Range("D63").Value2 = 0
$ws.Range("D64").Value2 = 0
$ws.Range("D65").Value2 = 0
$ws.Range("D66").Value2 = 11131000
$ws.Range("D67").ClearContents()
$ws.Range("D68").Value2 = 0
$ws.Range("D69").Value2 = 0
$ws.Range("D70").Value2 = 22000
$ws.Range("D71").Value2 = 0
$ws.Range("D72").Value2 = 3390000
$ws.Range("D73").Value2 = 0
$ws.Range("D74").Value2 = 0
$ws.Range("D75").Value2 = 0
$ws.Range("D76").Value2 = 4236000
$ws.Range("D77").Value2 = 0
$ws.Range("D80").Value2 = 43465
$ws.Range("D81").Value2 = 400000
$ws.Range("D82").ClearContents()
$ws.Range("D83").Value2 = 556000
$ws.Range("D84").Value2 = 0
$ws.Range("D85").Value2 = 0
$ws.Range("D86").Value2 = 0
$ws.Range("D87").Value2 = 0
$ws.Range("D88").Value2 = 0
$ws.Range("D89").Value2 = 1013000
$ws.Range("D90").ClearContents()
$ws.Range("D91").Value2 = -1538000
$ws.Range("D92").Value2 = 0
$ws.Range("D93").Value2 = 0
$ws.Range("D94").Value2 = -1531000
$ws.Range("D95").ClearContents()
$ws.Range("D96").Value2 = -51000
$ws.Range("D97").Value2 = 0
$ws.Range("D98").Value2 = 0
$ws.Range("D99").Value2 = 0
$ws.Range("D100").Value2 = 528000
$ws.Range("D101").Value2 = 0
$ws.Range("D102").Value2 = 10000
